$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "notes" header in column L
$ws.Range("L1").Value = "notes"

# Salary / salary-year data for ranks 26-40 (rows 27-41), continuing the
# existing J (salary) / K (salary-year) columns that were already
# populated through row 26.

# Row 27 - Baltimore, Maryland
$ws.Range("J27").Value = 185000
$ws.Range("K27").Value = "FY2019"

# Row 28 - Oklahoma City, Oklahoma
$ws.Range("J28").Value = 24000
$ws.Range("K28").Value = "FY2018"

# Row 29 - Louisville, Kentucky
$ws.Range("J29").Value = 126486.36
$ws.Range("K29").Value = 2019

# Row 30 - Portland, Oregon
$ws.Range("J30").Value = 143666
$ws.Range("K30").Value = "FY2018"

# Row 31 - (City commission)
$ws.Range("J31").Value = 144723.74
$ws.Range("K31").Value = 2018

# Row 32 - Las Vegas, Nevada
$ws.Range("J32").Value = 147335.76
$ws.Range("K32").Value = "FY2016"

# Row 33 - Milwaukee, Wisconsin
$ws.Range("J33").Value = 125000
$ws.Range("K33").Value = "FY2017"
$ws.Range("L33").Value = "Check more on this"

# Row 34 - Albuquerque, New Mexico
$ws.Range("J34").Value = 42000
$ws.Range("K34").Value = "FY2019"
$ws.Range("L34").Value = "Just confirmed via referendum "

# Row 35 - Tucson, Arizona
$ws.Range("J35").Value = 136900
$ws.Range("K35").Value = "FY2018"

# Row 36 - Fresno, California
$ws.Range("J36").Value = 130276
$ws.Range("K36").Value = 2018

# Row 37
$ws.Range("J37").Value = 129391
$ws.Range("K37").Value = 2018

# Row 38
$ws.Range("J38").Value = 141455
$ws.Range("K38").Value = "FY2019"

# Row 39
$ws.Range("J39").Value = 73000
$ws.Range("K39").Value = "FY2018"

# Row 40
$ws.Range("J40").Value = 30000
$ws.Range("K40").Value = "FY2018"

# Row 41
$ws.Range("J41").Value = 147500
$ws.Range("K41").Value = "FY2017"

# Reflect the cursor position the author ended up at when saving.
[void]$ws.Range("J42").Select()
